$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration_template")

# Insert two new blank rows at the top (rows 2 and 3), pushing the existing
# rows 2-18 down to rows 4-20.
$ws.Rows("2:3").Insert()

# Row 2 stays blank; row 3 gets the new "Unnamed: 0" label.
$ws.Range("A3").Value = "Unnamed: 0"

# Append the new rows 21-26 at the bottom of the list.
$ws.Range("A21").Value = "filter1"

$ws.Range("A22").Value = "filter2"

$ws.Range("A23").Value = "60m Wind Speed CorrWS"
$ws.Range("B23").Value = "corrWS_RSD_WS"

$ws.Range("A24").Value = "60m Wind Turbulence CorrWS"
$ws.Range("B24").Value = "corrWS_RSD_TI"

$ws.Range("A25").Value = "60m Wind Speed Corr TI"
$ws.Range("B25").Value = "corrTI_RSD_WS"

$ws.Range("A26").Value = "60m Wind Turbulence CorrTI"
$ws.Range("B26").Value = "corrTI_RSD_TI"

# Match the author's final selection.
$null = $ws.Range("B24").Select()
